$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 1472.5
$ws.Range("I18").Value = 1472.5
$ws.Range("K18").Value = 1472.5
$ws.Range("M18").Value = -1188.5
$ws.Range("H40").Value = 2805
$ws.Range("I40").Value = 4884
$ws.Range("J40").Value = 1765.5
$ws.Range("K40").Value = 4884
$ws.Range("L40").Value = 1765.5
$ws.Range("M40").Value = -4709
$ws.Range("N40").Value = -2115.5
$ws.Range("H53").Value = 732.8095
$ws.Range("J53").Value = 1480.125
$ws.Range("L53").Value = 1480.125
$ws.Range("N53").Value = -2754.125
$ws.Range("H103").Value = 397.15
$ws.Range("I103").Value = 406.7647
$ws.Range("J103").Value = 342.66666
$ws.Range("K103").Value = 1220.2941
$ws.Range("L103").Value = 1027.99998
$ws.Range("M103").Value = -634.2941000000001
$ws.Range("N103").Value = -2199.99998

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 5051512.5
$ws.Range("I2").Value = 6494516.5
$ws.Range("J2").Value = 999.5
$ws.Range("K2").Value = 6494516.5
$ws.Range("L2").Value = 999.5
$ws.Range("M2").Value = -6494403.5
$ws.Range("N2").Value = -1225.5
$ws.Range("H34").Value = 0
$ws.Range("I34").Value = 0
$ws.Range("K34").Value = 0
$ws.Range("M34").ClearContents()
$ws.Range("H116").Value = 5051512.5
$ws.Range("I116").Value = 6494516.5
$ws.Range("J116").Value = 999.5
$ws.Range("K116").Value = 6494516.5
$ws.Range("L116").Value = 999.5
$ws.Range("M116").Value = -6492222.5
$ws.Range("N116").Value = -5587.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 5051512.5
$ws.Range("I3").Value = 6494516.5
$ws.Range("J3").Value = 999.5
$ws.Range("K3").Value = 6494516.5
$ws.Range("L3").Value = 999.5
$ws.Range("M3").Value = -6494402.5
$ws.Range("N3").Value = -1227.5
$ws.Range("H35").Value = 0
$ws.Range("J35").Value = 0
$ws.Range("L35").Value = 0
$ws.Range("N35").ClearContents()
$ws.Range("H80").Value = 61460.176
$ws.Range("I80").Value = 521.2857
$ws.Range("J80").Value = 104117.4
$ws.Range("K80").Value = 521.2857
$ws.Range("L80").Value = 104117.4
$ws.Range("M80").Value = 476.7143
$ws.Range("N80").Value = -106113.4
$ws.Range("H83").Value = 61460.176
$ws.Range("I83").Value = 521.2857
$ws.Range("J83").Value = 104117.4
$ws.Range("K83").Value = 2606.4285
$ws.Range("L83").Value = 520587
$ws.Range("M83").Value = 2385.5715
$ws.Range("N83").Value = -530571
$ws.Range("H99").Value = 3100
$ws.Range("I99").Value = 2800
$ws.Range("K99").Value = 2800
$ws.Range("M99").Value = -1302
$ws.Range("I105").Value = 38463024
$ws.Range("J105").Value = 3044
$ws.Range("K105").Value = 38463024
$ws.Range("L105").Value = 3044
$ws.Range("M105").Value = -38461277
$ws.Range("N105").Value = -6538
$ws.Range("H117").Value = 94942
$ws.Range("J117").Value = 94942
$ws.Range("L117").Value = 94942
$ws.Range("N117").Value = -104120

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 4640.5835
$ws.Range("I58").Value = 2772
$ws.Range("K58").Value = 2772
$ws.Range("M58").Value = -2569
$ws.Range("H134").Value = 4456.579
$ws.Range("I134").Value = 3963.6365
$ws.Range("K134").Value = 11890.9095
$ws.Range("M134").Value = -9355.9095
$ws.Range("H136").Value = 4640.5835
$ws.Range("I136").Value = 2772
$ws.Range("K136").Value = 8316
$ws.Range("M136").Value = -5766

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 1581788.4
$ws.Range("I4").Value = 433657.06
$ws.Range("J4").Value = 3878051
$ws.Range("K4").Value = 1300971.18
$ws.Range("L4").Value = 11634153
$ws.Range("M4").Value = -1300859.18
$ws.Range("N4").Value = -11634377
$ws.Range("H62").Value = 2058.127
$ws.Range("I62").Value = 1229.7949
$ws.Range("J62").Value = 3404.1667
$ws.Range("K62").Value = 3689.384700000001
$ws.Range("L62").Value = 10212.5001
$ws.Range("M62").Value = -3003.384700000001
$ws.Range("N62").Value = -11584.5001
$ws.Range("H64").Value = 3530.2727
$ws.Range("J64").Value = 3483.3
$ws.Range("L64").Value = 10449.9
$ws.Range("N64").Value = -10989.9
$ws.Range("H65").Value = 2058.127
$ws.Range("I65").Value = 1229.7949
$ws.Range("J65").Value = 3404.1667
$ws.Range("K65").Value = 11068.1541
$ws.Range("L65").Value = 30637.5003
$ws.Range("M65").Value = -7636.154100000002
$ws.Range("N65").Value = -37501.5003
$ws.Range("H67").Value = 3530.2727
$ws.Range("J67").Value = 3483.3
$ws.Range("L67").Value = 10449.9
$ws.Range("N67").Value = -12321.9
$ws.Range("H107").Value = 2651.8
$ws.Range("I107").Value = 987.3333
$ws.Range("J107").Value = 3067.9167
$ws.Range("K107").Value = 2961.9999
$ws.Range("L107").Value = 9203.750100000001
$ws.Range("M107").Value = -1041.9999
$ws.Range("N107").Value = -13043.7501

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3939.0715
$ws.Range("I80").Value = 3048.7144
$ws.Range("J80").Value = 4829.4287
$ws.Range("K80").Value = 3048.7144
$ws.Range("L80").Value = 4829.4287
$ws.Range("M80").Value = -2050.7144
$ws.Range("N80").Value = -6825.4287
$ws.Range("H83").Value = 3939.0715
$ws.Range("I83").Value = 3048.7144
$ws.Range("J83").Value = 4829.4287
$ws.Range("K83").Value = 15243.572
$ws.Range("L83").Value = 24147.1435
$ws.Range("M83").Value = -10251.572
$ws.Range("N83").Value = -34131.14350000001
$ws.Range("H105").Value = 79999
$ws.Range("J105").Value = 79999
$ws.Range("L105").Value = 79999
$ws.Range("N105").Value = -86987
$ws.Range("H113").Value = 2431.8262
$ws.Range("I113").Value = 2474.8667
$ws.Range("J113").Value = 2351.125
$ws.Range("K113").Value = 2474.8667
$ws.Range("L113").Value = 2351.125
$ws.Range("M113").Value = -304.8667
$ws.Range("N113").Value = -6691.125
$ws.Range("H132").Value = 8992.909
$ws.Range("I132").Value = 6560.2856
$ws.Range("K132").Value = 19680.8568
$ws.Range("M132").Value = -17150.8568

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 777.1111
$ws.Range("I16").Value = 760.8461
$ws.Range("J16").Value = 1200
$ws.Range("K16").Value = 760.8461
$ws.Range("L16").Value = 1200
$ws.Range("M16").Value = -590.8461
$ws.Range("N16").Value = -1540
$ws.Range("H46").Value = 7649.087
$ws.Range("I46").Value = 1650
$ws.Range("J46").Value = 9315.5
$ws.Range("K46").Value = 1650
$ws.Range("L46").Value = 9315.5
$ws.Range("M46").Value = -1462
$ws.Range("N46").Value = -9691.5
$ws.Range("H133").Value = 109998.5
$ws.Range("J133").Value = 109998
$ws.Range("L133").Value = 109998
$ws.Range("N133").Value = -115058
$ws.Range("H136").Value = 2670905
$ws.Range("I136").Value = 4169289.8
$ws.Range("K136").Value = 12507869.4
$ws.Range("M136").Value = -12505319.4

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H47").Value = 18375
$ws.Range("J47").Value = 18375
$ws.Range("L47").Value = 18375
$ws.Range("N47").Value = -19519
$ws.Range("H116").Value = 92859
$ws.Range("J116").Value = 92859
$ws.Range("L116").Value = 92859
$ws.Range("N116").Value = -102037
$ws.Range("H126").Value = 6045.4814
$ws.Range("I126").Value = 5730.619
$ws.Range("K126").Value = 17191.857
$ws.Range("M126").Value = -14721.857
